# R1001_BOM.xlsx — read for version 1.1.2, crowdsupply build
#
# Summary of changes applied to the BOM sheet:
#  - C14/C9 100nF-cap group gains C112 (qty 2 -> 3)
#  - J1 (PCIEXPRESS-X1 connector) row removed entirely (no longer stuffed)
#  - J2 footprint corrected to Main:70553-0038
#  - New resistor line added: R120 R121 (100 kOhm, 0603, STUFF)
#  - R18 resistor group gains R4 (qty 1 -> 2)
#  - R4 and R21 removed from the big 10.0 kOhm resistor group (qty 13 -> 11)
#  - R20 (15.0 kOhm) row removed entirely
#  - U3 MCU part number bumped from EFM8BB10F2G to EFM8BB10F8G

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) C14 C9  -> C112 C14 C9  (row 5), quantity 2 -> 3
$ws.Range("A5").Value2 = "C112 C14 C9 "
$ws.Range("B5").Value2 = 3

# 2) Remove the J1 / PCIEXPRESS-X1 row (row 6) entirely
$ws.Rows.Item(6).EntireRow.Delete()

# 3) J2 footprint correction (J2 is now row 6 after the delete above)
$ws.Range("D6").Value2 = "Main:70553-0038"

# 4) Insert a new row for R120/R121 just above the R15 row.
#    R15 was row 10 before the J1 delete, so it is row 9 now -> insert at 9.
$ws.Rows.Item(9).EntireRow.Insert()

$ws.Range("A9").Value2 = "R120 R121 "
$ws.Range("B9").Value2 = 2
$ws.Range("C9").Value2 = "RC0603FR-07100KL"
$ws.Range("D9").Value2 = "Main:RESC0603"
$ws.Range("F9").Value2 = "Yageo"
$ws.Range("G9").Value2 = "RC0603FR-07100KL"
$ws.Range("H9").Value2 = "STUFF"
$ws.Range("J9").Value2 = "±1%"
$ws.Range("M9").Value2 = "100 kΩ"
$ws.Range("N9").Value2 = "100 mW"

# 5) R18 group gains R4: "R18 " -> "R18 R4 ", quantity 1 -> 2 (now row 11)
$ws.Range("A11").Value2 = "R18 R4 "
$ws.Range("B11").Value2 = 2

# 6) Big 10.0 kOhm group loses R21 and R4 (moved to R18 group above); quantity 13 -> 11 (now row 12)
$ws.Range("A12").Value2 = "R1 R10 R11 R12 R13 R14 R16 R17 R2 R6 R7 "
$ws.Range("B12").Value2 = 11

# 7) Remove the R20 (15.0 kOhm) row entirely (now row 13)
$ws.Rows.Item(13).EntireRow.Delete()

# 8) U3 part-number bump: EFM8BB10F2G -> EFM8BB10F8G (now row 18)
$ws.Range("C18").Value2 = "EFM8BB10F8G"
$ws.Range("G18").Value2 = "EFM8BB10F8G-A-QFN20R"

# 9) Column A narrows slightly now that the longest reference strings are gone
$ws.Columns.Item(1).ColumnWidth = 38.14

# 10) Restore the author's last-saved selection
$ws.Range("D7").Select()
